$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 472; this shifts rows 472:500 down to 473:501
$ws.Rows.Item(472).Insert()

# Populate the new row 472 with the new weekly record
$ws.Range("A472").Value = 10
$ws.Range("B472").Value = "Vega Modelo de Temuco"
$ws.Range("C472").Value = "La Araucanía"
$ws.Range("D472").Value = 44509
$ws.Range("D472").NumberFormat = $ws.Range("D473").NumberFormat
$ws.Range("E472").Value = 9
$ws.Range("F472").Value = 100112021
$ws.Range("G472").Value = "Ají"
$ws.Range("H472").Value = "Americana (o)"
$ws.Range("I472").Value = "Primera"
$ws.Range("J472").Value = 20
$ws.Range("K472").Value = 45000
$ws.Range("L472").Value = 45000
$ws.Range("M472").Value = 45000
$ws.Range("N472").Value = "$/caja 25 kilos"
$ws.Range("O472").Value = "Provincia de Limarí"
$ws.Range("P472").Value = 1800
$ws.Range("Q472").Value = 25
$ws.Range("R472").Value = "Hortaliza"
